$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 4399.25
$ws.Range("I12").Value = 4998.095
$ws.Range("J12").Value = 207.33333
$ws.Range("K12").Value = 4998.095
$ws.Range("L12").Value = 207.33333
$ws.Range("M12").Value = -4828.095
$ws.Range("N12").Value = -547.3333299999999
$ws.Range("H19").Value = 1104
$ws.Range("I19").Value = 1040.6666
$ws.Range("K19").Value = 1040.6666
$ws.Range("M19").Value = -865.6666
$ws.Range("H28").Value = 823.2
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H29").Value = 375.77777
$ws.Range("H38").Value = 497.83334
$ws.Range("I38").Value = 497.83334
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 1493.50002
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -1121.50002
$ws.Range("N38").ClearContents()
$ws.Range("H51").Value = 4682.6665
$ws.Range("I51").Value = 4682.6665
$ws.Range("K51").Value = 4682.6665
$ws.Range("M51").Value = -4198.6665
$ws.Range("H58").Value = 71473.71000000001
$ws.Range("I58").Value = 78.75
$ws.Range("J58").Value = 166667
$ws.Range("K58").Value = 236.25
$ws.Range("L58").Value = 500001
$ws.Range("M58").Value = -86.25
$ws.Range("N58").Value = -500301
$ws.Range("H132").Value = 3868.158
$ws.Range("I132").Value = 3868.158
$ws.Range("K132").Value = 11604.474
$ws.Range("M132").Value = -9074.474
$ws.Range("H137").Value = 2231.5625
$ws.Range("J137").Value = 3912.7144
$ws.Range("L137").Value = 11738.1432
$ws.Range("N137").Value = -16838.1432

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 62503150
$ws.Range("J61").Value = 7072.25
$ws.Range("L61").Value = 7072.25
$ws.Range("N61").Value = -7496.25
$ws.Range("H74").Value = 23257422
$ws.Range("J74").Value = 1961.3
$ws.Range("L74").Value = 1961.3
$ws.Range("N74").Value = -3709.3
$ws.Range("H77").Value = 23257422
$ws.Range("J77").Value = 1961.3
$ws.Range("L77").Value = 9806.5
$ws.Range("N77").Value = -18542.5
$ws.Range("H110").Value = 44294.39
$ws.Range("I110").Value = 48440.145
$ws.Range("K110").Value = 48440.145
$ws.Range("M110").Value = -46395.145
$ws.Range("H136").Value = 62503150
$ws.Range("J136").Value = 7072.25
$ws.Range("L136").Value = 21216.75
$ws.Range("N136").Value = -26316.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 6266.5
$ws.Range("I102").Value = 6266.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 6266.5
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -3021.5
$ws.Range("N102").ClearContents()
$ws.Range("H105").Value = 3324.625
$ws.Range("I105").Value = 3044.9092
$ws.Range("J105").Value = 3940
$ws.Range("K105").Value = 3044.9092
$ws.Range("L105").Value = 3940
$ws.Range("M105").Value = -1297.9092
$ws.Range("N105").Value = -7434
$ws.Range("H134").Value = 25001404
$ws.Range("I134").Value = 25001404
$ws.Range("K134").Value = 75004212
$ws.Range("M134").Value = -75001677

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6052.2144
$ws.Range("I31").Value = 8578.333000000001
$ws.Range("K31").Value = 8578.333000000001
$ws.Range("M31").Value = -8283.333000000001
$ws.Range("H34").Value = 6052.2144
$ws.Range("I34").Value = 8578.333000000001
$ws.Range("K34").Value = 8578.333000000001
$ws.Range("M34").Value = -8376.333000000001
$ws.Range("H86").Value = 10075.091
$ws.Range("I86").Value = 7457.125
$ws.Range("J86").Value = 11571.071
$ws.Range("K86").Value = 7457.125
$ws.Range("L86").Value = 11571.071
$ws.Range("M86").Value = -6334.125
$ws.Range("N86").Value = -13817.071
$ws.Range("H89").Value = 10075.091
$ws.Range("I89").Value = 7457.125
$ws.Range("J89").Value = 11571.071
$ws.Range("K89").Value = 37285.625
$ws.Range("L89").Value = 57855.355
$ws.Range("M89").Value = -31669.625
$ws.Range("N89").Value = -69087.355
$ws.Range("H132").Value = 27029318
$ws.Range("I132").Value = 28573764
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 85721292
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -85718762
$ws.Range("N132").Value = -9560
$ws.Range("H141").Value = 88243.125
$ws.Range("I141").Value = 93934
$ws.Range("K141").Value = 93934
$ws.Range("M141").Value = -88754

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 423.7143
$ws.Range("I92").Value = 430
$ws.Range("J92").Value = 419
$ws.Range("K92").Value = 1290
$ws.Range("L92").Value = 1257
$ws.Range("M92").Value = -42
$ws.Range("N92").Value = -3753
$ws.Range("I132").Value = 60000
$ws.Range("J132").Value = 2319.4
$ws.Range("K132").Value = 540000
$ws.Range("L132").Value = 20874.6
$ws.Range("M132").Value = -537470
$ws.Range("N132").Value = -25934.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 6900
$ws.Range("J26").Value = 6900
$ws.Range("L26").Value = 6900
$ws.Range("N26").Value = -7460
$ws.Range("H50").Value = 6900
$ws.Range("J50").Value = 6900
$ws.Range("L50").Value = 6900
$ws.Range("N50").Value = -7896
$ws.Range("H64").Value = 70000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 70000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 70000
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -70496
$ws.Range("H67").Value = 70000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 70000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 70000
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -71716
$ws.Range("H107").Value = 1078.1154
$ws.Range("I107").Value = 736.2632
$ws.Range("J107").Value = 2006
$ws.Range("K107").Value = 736.2632
$ws.Range("L107").Value = 2006
$ws.Range("M107").Value = 1183.7368
$ws.Range("N107").Value = -5846

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 21875742
$ws.Range("I100").Value = 21875742
$ws.Range("K100").Value = 21875742
$ws.Range("M100").Value = -21875201
$ws.Range("H130").Value = 96795
$ws.Range("J130").Value = 96795
$ws.Range("L130").Value = 96795
$ws.Range("N130").Value = -106835

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1749.8334
$ws.Range("I81").Value = 1699.75
$ws.Range("J81").Value = 1850
$ws.Range("K81").Value = 3399.5
$ws.Range("L81").Value = 3700
$ws.Range("M81").Value = -2338.5
$ws.Range("N81").Value = -5822
$ws.Range("H84").Value = 1749.8334
$ws.Range("I84").Value = 1699.75
$ws.Range("J84").Value = 1850
$ws.Range("K84").Value = 16997.5
$ws.Range("L84").Value = 18500
$ws.Range("M84").Value = -11693.5
$ws.Range("N84").Value = -29108
$ws.Range("H136").Value = 26317092
$ws.Range("I136").Value = 29412956
$ws.Range("K136").Value = 88238868
$ws.Range("M136").Value = -88236318

Write-Output "Applied Spriggan_Profits market-data refresh."